# "finish student tests for hw3"
#
# - Rename the display_name_* columns (B/C) to is_list_ordered_*
# - Insert 8 new columns (new D..K) holding 4 more score/review pairs
#   (picking_right_pair, picking_wrong_pair, picking_zero_pair,
#   picking_only_one_element); this pushes the old final_score column
#   from D to L
# - Fill in the score/review values for both students and update
#   final_score

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new columns -----------------------------------------
# Original layout: A=student_ids, B=display_name_score, C=display_name_review,
#                   D=final_score
# Insert 8 columns before the existing D column (final_score), which then
# slides from D to L, picking up the existing header's bold/border style on
# the way.
$ws.Range("D1:K1").EntireColumn.Insert()

# --- Header row (row 1) ------------------------------------------------------
$ws.Range("B1").Value = "is_list_ordered_score"
$ws.Range("C1").Value = "is_list_ordered_review"
$ws.Range("D1").Value = "picking_right_pair_score"
$ws.Range("E1").Value = "picking_right_pair_review"
$ws.Range("F1").Value = "picking_wrong_pair_score"
$ws.Range("G1").Value = "picking_wrong_pair_review"
$ws.Range("H1").Value = "picking_zero_pair_score"
$ws.Range("I1").Value = "picking_zero_pair_review"
$ws.Range("J1").Value = "picking_only_one_element_score"
$ws.Range("K1").Value = "picking_only_one_element_review"
$ws.Range("L1").Value = "final_score"

# --- Row 2 (student 204897687) ----------------------------------------------
$ws.Range("B2").Value = "20.0 / 20.0"
$ws.Range("D2").Value = "4.0 / 20.0"
$ws.Range("E2").Value = "bruce_wayne_and_wayne_enterprises: failed`nclark_kent_and_daily_planet: failed`npeter_parker_and_daily_bugle: failed`nwillie_wonka_and_chocolate_factory: failed"
$ws.Range("F2").Value = "0.0 / 20.0"
$ws.Range("G2").Value = "bruce_wayne_and_chocolate_factory: failed"
$ws.Range("H2").Value = "20.0 / 20.0"
$ws.Range("J2").Value = "20.0 / 20.0"
$ws.Range("L2").Value = "64.0/100.0"

# --- Row 3 (student 308418367) ----------------------------------------------
$ws.Range("B3").Value = "20.0 / 20.0"
$ws.Range("D3").Value = "20.0 / 20.0"
$ws.Range("F3").Value = "20.0 / 20.0"
$ws.Range("H3").Value = "20.0 / 20.0"
$ws.Range("J3").Value = "20.0 / 20.0"
$ws.Range("L3").Value = "100.0/100.0"

# --- Empty (but present) review cells ---------------------------------------
# C2 held the old "display_name_review" failure text; drop it since the new
# is_list_ordered_review result is blank (no failures) for this student.
$ws.Range("C2").ClearContents()

# C3 already existed as an empty cell in the original workbook; reuse it as a
# template (paste formats only, which is enough to materialize the cell
# without assigning it a value) so the blank cells (C2, I2, K2, E3, G3, I3,
# K3) stay present-but-empty like the rest of the "review" columns, instead
# of disappearing the way assigning "" / ClearContents alone would leave them.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
